$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medicos")

# Correcciones de ortografia (acentos) en los encabezados de columna D
$ws.Range("D3").Value = "Código Postal"
$ws.Range("D11").Value = "Número Interior"
$ws.Range("D9").Value = "Número Exterior"
$ws.Range("D21").Value = "Teléfono"

# Actualizar la seleccion / vista de la hoja (topLeftCell="A13", selection D21)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D21").Select()
